$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily 100-error-count rows appended below the existing data (rows 80-81)
$ws.Range("A80").Value = 46045
$ws.Range("B80").Value = 695
$ws.Range("C80").Value = 637
$ws.Range("D80").Value = 58
$ws.Range("A80").NumberFormat = "d-mmm-yy"

$ws.Range("A81").Value = 46048
$ws.Range("B81").Value = 1219
$ws.Range("C81").Value = 1201
$ws.Range("D81").Value = 18
$ws.Range("A81").NumberFormat = "d-mmm-yy"

# Match the saved view/selection state: scrolled down, new last row selected
$ws.Range("A68").Select()
$ws.Range("A81:D81").Select()
